$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valentin")

# Mark rows 17, 27, 38 as reserved/bought
$ws.Range("E17").Value = "Y"
$ws.Range("E27").Value = "Y"
$ws.Range("E38").Value = "Y"

# Add new wishlist item in row 40
$ws.Range("A40").Value = "Pulp Fiction [4K Ultra HD]"
$ws.Range("B40").Value = "https://m.media-amazon.com/images/I/61v0QDTS-kL._SX522_.jpg"
$ws.Range("C40").Value = "https://www.amazon.de/-/en/Pulp-Fiction-Limited-Collectors-Blu-ray/dp/B0DJDGYFTG?dib=eyJ2IjoiMSJ9.klW5y3HEDcmd7XsJ5SEDjsF3dIxtdAK9F9DnVk-09K-jG4H-CeSBPlOXcvVa7BzhIMirUVtfM90L_viZ6GLyyYeS8tLnSm4nWX4dICrV5ZAnTl50gG28XiIWyUvtdcU1WrtcuAoHHevKzxqu_iNKyy_HFqWlLiYxdnko2L-Iz6EGEeXPYOUjZfeScCFBTYAs6OTdF5RE9LRwK5ITWgS5Ag8eskzKt0ZdLhkH_pjbxXI.R25ziHAmaTXuk5wZzqxoNjqm7stqNzTBEkKJ00SARqo&dib_tag=se&keywords=Pulp+Fiction+30th+Anniversary+Collector%27s+Edition+Amazon+Exclusive+%5B4K+UHD+Blu-Ray+Digital+Copy%5D&linkCode=gg3&linkId=3ddeaf77413169d18b3d6c726143490d&nsdOptOutParam=true&qid=1735808493&sr=8-1&ufe=INHOUSE_INSTALLMENTS%3ADE_IHI_3M_AUTOMATED"
$ws.Range("D40").Value = "45 EUR"

$ws.Range("C46").Select()
